$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new "sex" value (re-ran anonymisation: swap Female/Male for specific rows)
$rowValues = @{
    4 = "Female"
    10 = "Male"
    13 = "Male"
    26 = "Female"
    27 = "Female"
    28 = "Male"
    37 = "Male"
    39 = "Female"
    50 = "Male"
    51 = "Male"
    53 = "Female"
    58 = "Male"
    60 = "Male"
    62 = "Male"
    63 = "Female"
    71 = "Male"
    72 = "Male"
    78 = "Female"
    80 = "Male"
    81 = "Male"
    82 = "Male"
    90 = "Male"
    92 = "Male"
    94 = "Male"
    102 = "Male"
    105 = "Female"
    108 = "Female"
    111 = "Female"
    115 = "Male"
    117 = "Male"
    120 = "Male"
    121 = "Male"
    128 = "Male"
    139 = "Male"
    145 = "Male"
    148 = "Female"
    160 = "Male"
    162 = "Male"
    165 = "Male"
    174 = "Female"
    177 = "Male"
    178 = "Female"
    182 = "Female"
    185 = "Male"
    192 = "Female"
    195 = "Female"
    198 = "Male"
}

foreach ($row in $rowValues.Keys) {
    $ws.Cells.Item($row, 1).Value = $rowValues[$row]
}
